$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data rows (do the copies before any values are overwritten) ---

# New column H ("season_y") duplicates the year column's text values
$ws.Range("A2:A6").Copy($ws.Range("H2:H6"))

# The existing pts_per_game values (currently column F) move to the new column G
$ws.Range("F2:F6").Copy($ws.Range("G2:G6"))

# Populate the new pf_per_game values into column F
$ws.Range("F2").Value2 = 2.1
$ws.Range("F3").Value2 = 3.3
$ws.Range("F4").Value2 = 2.9
$ws.Range("F5").Value2 = 3.4
$ws.Range("F6").Value2 = 2.2

# --- Header row ---

# Give the two new header cells (G1, H1) the same formatting as the other
# header cells (bold, centered, bordered) before setting their text.
$ws.Range("A1").Copy($ws.Range("G1"))
$ws.Range("A1").Copy($ws.Range("H1"))

# Rename the "season" header to "year"
$ws.Range("A1").Value2 = "year"

# F1 held "pts_per_game"; it's now the header for the new pf_per_game column
$ws.Range("F1").Value2 = "pf_per_game"

# G1 is the new header for the (moved) pts_per_game column
$ws.Range("G1").Value2 = "pts_per_game"

# H1 is the new header for the season_y column
$ws.Range("H1").Value2 = "season_y"
